# Daily attendance processing - 2026-01-15 06:09:23
# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in the "Recorded By" column (column G) wherever it appears exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

$rng = $ws.Range("G1:G259")
$found = $rng.Find($target)
if ($found -ne $null) {
    $firstAddress = $found.Address()
    $continue = $true
    while ($continue) {
        $found.Value2 = $replacement
        $found = $rng.FindNext($found)
        if ($found -eq $null -or $found.Address() -eq $firstAddress) {
            $continue = $false
        }
    }
}
